# Refresh the crypto price/volume snapshot values (columns D and E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '82.137.13'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +3.22%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.196.11'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.16'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +5.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '627.45'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.26%  '
$ws.Range('E7').Value = '  +21.62%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.587'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '3.194.84'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.590'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000260'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +13.65%  '
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.33'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.779.90'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '31.79'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '81.905.87'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.14%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.191.21'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('E19').Value = '  +5.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.07'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '436.05'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.88%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.97'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.12'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('E24').Value = '  +5.88%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.29'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +11.17%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.349.27'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '76.78'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.97'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.54%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '590.02'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +12.71%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.05'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.156'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +8.31%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.01'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('E37').Value = '  +15.45%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '22.84'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.14'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +11.02%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.409'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.05'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +13.90%  '
$ws.Range('E43').Value = '  +23.02%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '20.80'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '160.71'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '188.33'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '44.59'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.42%  '
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('E50').Value = '  -5.60%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '26.30'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.94%  '
